# Fruta / hortaliza, semanal
# Insert two new weekly price-report rows above the existing "Vega Modelo
# de Temuco - Poroto verde" records (pushing the former rows 110-116 down
# to 112-118), then populate the two new rows with this week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing rows 110:116 down by two rows.
$ws.Rows("110:111").Insert()

# New row 110
$ws.Range("A110").Value = 10
$ws.Range("B110").Value = "Vega Modelo de Temuco"
$ws.Range("C110").Value = "La Araucanía"
$ws.Range("D110").Value = 44610
$ws.Range("E110").Value = 9
$ws.Range("F110").Value = 100112031
$ws.Range("G110").Value = "Poroto verde"
$ws.Range("H110").Value = "Brío"
$ws.Range("I110").Value = "Primera"
$ws.Range("J110").Value = 50
$ws.Range("K110").Value = 1200
$ws.Range("L110").Value = 1200
$ws.Range("M110").Value = 1200
$ws.Range("N110").Value = "$/kilo"
$ws.Range("O110").Value = "Región de La Araucanía"
$ws.Range("P110").Value = 1200
$ws.Range("Q110").Value = 1
$ws.Range("R110").Value = "Hortaliza"

# New row 111
$ws.Range("A111").Value = 10
$ws.Range("B111").Value = "Vega Modelo de Temuco"
$ws.Range("C111").Value = "La Araucanía"
$ws.Range("D111").Value = 44610
$ws.Range("E111").Value = 9
$ws.Range("F111").Value = 100112031
$ws.Range("G111").Value = "Poroto verde"
$ws.Range("H111").Value = "Sin especificar"
$ws.Range("I111").Value = "Primera"
$ws.Range("J111").Value = 80
$ws.Range("K111").Value = 1200
$ws.Range("L111").Value = 1200
$ws.Range("M111").Value = 1200
$ws.Range("N111").Value = "$/kilo"
$ws.Range("O111").Value = "Región de La Araucanía"
$ws.Range("P111").Value = 1200
$ws.Range("Q111").Value = 1
$ws.Range("R111").Value = "Hortaliza"
